$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "Price" (D) and "Volume(1h)" (E) columns with the latest scrape.
# These columns hold plain text (e.g. "69.417.25", "  +1.85%  "), not numbers,
# so every value is written with a leading apostrophe. Excel/COM treats a
# leading ' as a 'force text' marker (it is stripped from the stored value) -
# without it, number-looking strings such as "580.40" would be auto-converted
# to the numeric value 580.4 instead of staying the literal text "580.40".
function Set-TextValue($range, [string]$text) {
    $range.Value = "'" + $text
}

Set-TextValue $ws.Range("D2") '69.417.25'
Set-TextValue $ws.Range("E2") '  +1.85%  '
Set-TextValue $ws.Range("D3") '3.378.21'
Set-TextValue $ws.Range("E3") '  +1.13%  '
Set-TextValue $ws.Range("E4") '  +0.09%  '
Set-TextValue $ws.Range("D5") '580.40'
Set-TextValue $ws.Range("E5") '  -0.48%  '
Set-TextValue $ws.Range("D6") '178.19'
Set-TextValue $ws.Range("E6") '  +0.56%  '
Set-TextValue $ws.Range("E7") '  +0.02%  '
Set-TextValue $ws.Range("E8") '  +0.49%  '
Set-TextValue $ws.Range("D9") '0.197'
Set-TextValue $ws.Range("E9") '  +8.14%  '
Set-TextValue $ws.Range("D10") '0.586'
Set-TextValue $ws.Range("E10") '  +0.78%  '
Set-TextValue $ws.Range("D11") '48.41'
Set-TextValue $ws.Range("E11") '  +0.79%  '
Set-TextValue $ws.Range("D12") '0.0000283'
Set-TextValue $ws.Range("E12") '  +3.81%  '
Set-TextValue $ws.Range("D13") '688.41'
Set-TextValue $ws.Range("E13") '  -0.81%  '
Set-TextValue $ws.Range("D14") '8.59'
Set-TextValue $ws.Range("E14") '  +2.08%  '
Set-TextValue $ws.Range("D15") '3.925.42'
Set-TextValue $ws.Range("E15") '  +0.85%  '
Set-TextValue $ws.Range("D16") '69.521.24'
Set-TextValue $ws.Range("E16") '  +1.90%  '
Set-TextValue $ws.Range("D18") '3.381.15'
Set-TextValue $ws.Range("E18") '  +0.23%  '
Set-TextValue $ws.Range("D19") '17.73'
Set-TextValue $ws.Range("E19") '  +1.73%  '
Set-TextValue $ws.Range("D20") '11.24'
Set-TextValue $ws.Range("E20") '  +0.68%  '
Set-TextValue $ws.Range("D21") '0.908'
Set-TextValue $ws.Range("E21") '  +1.52%  '
Set-TextValue $ws.Range("D22") '17.22'
Set-TextValue $ws.Range("E22") '  +1.31%  '
Set-TextValue $ws.Range("E23") '  -2.21%  '
Set-TextValue $ws.Range("D24") '101.35'
Set-TextValue $ws.Range("E24") '  +1.29%  '
Set-TextValue $ws.Range("E25") '  -0.77%  '
Set-TextValue $ws.Range("E26") '  -0.37%  '
Set-TextValue $ws.Range("D27") '9.68'
Set-TextValue $ws.Range("E27") '  +1.87%  '
Set-TextValue $ws.Range("D28") '33.50'
Set-TextValue $ws.Range("E28") '  +1.75%  '
Set-TextValue $ws.Range("D29") '8.71'
Set-TextValue $ws.Range("E29") '  +2.54%  '
Set-TextValue $ws.Range("D30") '6.90'
Set-TextValue $ws.Range("E30") '  -0.33%  '
Set-TextValue $ws.Range("D31") '3.83'
Set-TextValue $ws.Range("E31") '  +16.85%  '
Set-TextValue $ws.Range("D32") '11.03'
Set-TextValue $ws.Range("E32") '  +0.00%  '
Set-TextValue $ws.Range("D33") '553.57'
Set-TextValue $ws.Range("E33") '  -1.96%  '
Set-TextValue $ws.Range("E34") '  +0.09%  '
Set-TextValue $ws.Range("D35") '57.80'
Set-TextValue $ws.Range("E35") '  +0.47%  '
Set-TextValue $ws.Range("D36") '0.999'
Set-TextValue $ws.Range("E36") '  -0.03%  '
Set-TextValue $ws.Range("D37") '3.602.38'
Set-TextValue $ws.Range("E37") '  -2.35%  '
Set-TextValue $ws.Range("E38") '  +2.87%  '
Set-TextValue $ws.Range("D39") '35.29'
Set-TextValue $ws.Range("E39") '  +1.63%  '
Set-TextValue $ws.Range("D40") '0.0₃0726'
Set-TextValue $ws.Range("E40") '  +8.19%  '
Set-TextValue $ws.Range("D41") '3.30'
Set-TextValue $ws.Range("E41") '  +4.22%  '
Set-TextValue $ws.Range("D42") '2.72'
Set-TextValue $ws.Range("E42") '  +4.23%  '
Set-TextValue $ws.Range("D43") '0.0424'
Set-TextValue $ws.Range("E43") '  +2.63%  '
Set-TextValue $ws.Range("D44") '0.335'
Set-TextValue $ws.Range("E44") '  +0.12%  '
Set-TextValue $ws.Range("D45") '2.65'
Set-TextValue $ws.Range("E45") '  +0.20%  '
Set-TextValue $ws.Range("E46") '  +0.27%  '
Set-TextValue $ws.Range("E47") '  -0.18%  '
Set-TextValue $ws.Range("E48") '  +3.61%  '
Set-TextValue $ws.Range("D49") '128.95'
Set-TextValue $ws.Range("E49") '  -1.54%  '
Set-TextValue $ws.Range("D50") '2.57'
Set-TextValue $ws.Range("E50") '  -0.14%  '
Set-TextValue $ws.Range("D51") '7.40'
Set-TextValue $ws.Range("E51") '  -0.52%  '
